# Auto-generated script to update cryptos worksheet values
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('D2').Value = '30.752.20'
$ws.Range('E2').Value = '  +2.12%  '
$ws.Range('D3').Value = '1.687.47'
$ws.Range('E3').Value = '  +2.62%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = "'220.97"
$ws.Range('E5').Value = '  +2.36%  '
$ws.Range('E6').Value = '  -0.06%  '
$ws.Range('D7').Value = "'0.999"
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').Value = "'30.79"
$ws.Range('E8').Value = '  +4.91%  '
$ws.Range('E9').Value = '  +1.61%  '
$ws.Range('D10').Value = "'0.0625"
$ws.Range('E10').Value = '  +1.72%  '
$ws.Range('D11').Value = "'0.0906"
$ws.Range('E11').Value = '  -1.07%  '
$ws.Range('D12').Value = '1.932.89'
$ws.Range('E12').Value = '  +2.84%  '
$ws.Range('D13').Value = "'10.56"
$ws.Range('E13').Value = '  +11.46%  '
$ws.Range('D14').Value = "'0.629"
$ws.Range('E14').Value = '  +9.55%  '
$ws.Range('D15').Value = '1.704.59'
$ws.Range('E15').Value = '  +3.62%  '
$ws.Range('E16').Value = '  +1.81%  '
$ws.Range('D17').Value = '30.763.42'
$ws.Range('E17').Value = '  +2.09%  '
$ws.Range('D18').Value = "'66.56"
$ws.Range('E18').Value = '  +2.44%  '
$ws.Range('D19').Value = "'247.45"
$ws.Range('E19').Value = '  -0.23%  '
$ws.Range('D20').Value = '0.0₃0715'
$ws.Range('E20').Value = '  +0.55%  '
$ws.Range('D21').Value = "'0.999"
$ws.Range('E21').Value = '  -0.01%  '
$ws.Range('D22').Value = "'10.25"
$ws.Range('E22').Value = '  +2.49%  '
$ws.Range('D23').Value = "'4.28"
$ws.Range('E23').Value = '  +1.70%  '
$ws.Range('E24').Value = '  +0.44%  '
$ws.Range('D25').Value = "'157.27"
$ws.Range('E25').Value = '  -1.21%  '
$ws.Range('D26').Value = "'15.88"
$ws.Range('E26').Value = '  +0.56%  '
$ws.Range('E27').Value = '  -0.23%  '
$ws.Range('E28').Value = '  +0.52%  '
$ws.Range('E29').Value = '  -0.03%  '
$ws.Range('E30').Value = '  +1.16%  '
$ws.Range('E31').Value = '  +0.34%  '
$ws.Range('E32').Value = '  +1.61%  '
$ws.Range('D33').Value = '1.513.79'
$ws.Range('E33').Value = '  +5.10%  '
$ws.Range('D34').Value = "'3.29"
$ws.Range('E34').Value = '  +2.17%  '
$ws.Range('E35').Value = '  +4.09%  '
$ws.Range('E36').Value = '  -0.83%  '
$ws.Range('D37').Value = "'83.62"
$ws.Range('E37').Value = '  +5.85%  '
$ws.Range('E38').Value = '  +4.14%  '
$ws.Range('D39').Value = "'2.74"
$ws.Range('E39').Value = '  -4.80%  '
$ws.Range('E40').Value = '  +4.21%  '
$ws.Range('E41').Value = '  +1.49%  '
$ws.Range('D42').Value = "'0.848"
$ws.Range('E42').Value = '  -0.11%  '
$ws.Range('D43').Value = "'2.00"
$ws.Range('E43').Value = '  -1.54%  '
$ws.Range('D44').Value = "'0.0501"
$ws.Range('E44').Value = '  +0.60%  '
$ws.Range('E45').Value = '  -1.08%  '
$ws.Range('D46').Value = "'0.999"
$ws.Range('E46').Value = '  -0.05%  '
$ws.Range('D47').Value = "'51.89"
$ws.Range('E47').Value = '  -6.72%  '
$ws.Range('D48').Value = '1.825.65'
$ws.Range('E48').Value = '  +2.16%  '
$ws.Range('E49').Value = '  +0.75%  '
$ws.Range('D50').Value = "'94.83"
$ws.Range('E50').Value = '  +4.72%  '
$ws.Range('E51').Value = '  +0.82%  '
